# Rotate the block of rows 7-13 (A:C) up by one row:
# Zion Williamson (currently row 7) moves down to row 13,
# and the six rows below it (Anthony Edwards ... Nicolas Claxton) each shift up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A7:C13")
$values = $srcRange.Value()

$rowCount = $values.GetUpperBound(0)   # 7 (1-based upper bound)
$colCount = $values.GetUpperBound(1)   # 3 (1-based upper bound)

$newValues = New-Object 'object[,]' $rowCount, $colCount

for ($c = 1; $c -le $colCount; $c++) {
    for ($r = 1; $r -lt $rowCount; $r++) {
        $newValues[($r - 1), ($c - 1)] = $values[($r + 1), $c]
    }
    $newValues[($rowCount - 1), ($c - 1)] = $values[1, $c]
}

$ws.Range("A7:C13").Value = $newValues
